$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.347.67'
$ws.Range("E2").Value = '  +0.10%  '

$ws.Range("D3").Value = '1.931.74'
$ws.Range("E3").Value = '  +0.01%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.31%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7500'
$ws.Range("E5").Value = '  +5.17%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '244.14'
$ws.Range("E6").Value = '  -3.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  +0.15%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3164'
$ws.Range("E8").Value = '  -3.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '27.44'
$ws.Range("E9").Value = '  -0.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06986'
$ws.Range("E10").Value = '  -3.04%  '

$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08016'
$ws.Range("E11").Value = '  -0.96%  '

$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7740'
$ws.Range("E12").Value = '  -3.16%  '

$ws.Range("D13").Value = '1.932.37'
$ws.Range("E13").Value = '  +0.02%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.326'
$ws.Range("E14").Value = '  -1.89%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.25'
$ws.Range("E15").Value = '  -0.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.33'
$ws.Range("E16").Value = '  -3.34%  '

$ws.Range("D17").Value = '30.345.19'
$ws.Range("E17").Value = '  +0.07%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '252.74'
$ws.Range("E18").Value = '  +0.28%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007913'
$ws.Range("E19").Value = '  -2.53%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.743'
$ws.Range("E20").Value = '  -0.80%  '

$ws.Range("D21").Value = '2.184.80'
$ws.Range("E21").Value = '  +0.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.11%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  +0.27%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.646'
$ws.Range("E24").Value = '  -4.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.455'
$ws.Range("E25").Value = '  -2.63%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.90'
$ws.Range("E26").Value = '  +0.68%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.97'
$ws.Range("E27").Value = '  -1.36%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1322'
$ws.Range("E28").Value = '  +2.97%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.185'
$ws.Range("E29").Value = '  -5.76%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.373'
$ws.Range("E30").Value = '  +0.97%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.514'

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.386'
$ws.Range("E32").Value = '  -1.06%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.107'
$ws.Range("E33").Value = '  -2.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05128'
$ws.Range("E34").Value = '  -1.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.278'
$ws.Range("E35").Value = '  +1.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7470'
$ws.Range("E36").Value = '  -0.14%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.778'
$ws.Range("E37").Value = '  +0.55%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01951'
$ws.Range("E38").Value = '  -0.71%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.802'
$ws.Range("E39").Value = '  +0.11%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.425'
$ws.Range("E40").Value = '  +0.05%  '

$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '77.14'
$ws.Range("E41").Value = '  -2.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4441'
$ws.Range("E42").Value = '  -1.84%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.955'
$ws.Range("E43").Value = '  -3.57%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.002'
$ws.Range("E44").Value = '  +0.19%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8336'
$ws.Range("E45").Value = '  -0.99%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.38'
$ws.Range("E46").Value = '  -1.39%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.731'
$ws.Range("E47").Value = '  -1.22%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.457'
$ws.Range("E48").Value = '  +0.40%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '37.35'
$ws.Range("E49").Value = '  +1.71%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '973.27'
$ws.Range("E50").Value = '  +9.92%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06048'
$ws.Range("E51").Value = '  -0.71%  '
